$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 9-11 labels first (this matches the order new shared strings were
#     introduced by the original author): rollingWindow, lag_order, forecast_horizon ---
$ws.Range("A11").Value = "rollingWindow"
$ws.Range("A9").Value = "lag_order"
$ws.Range("A10").Value = "forecast_horizon"

# --- Add header note for column C ---
$ws.Range("C1").Value = "NOTES"

# --- Notes column content for existing rows 2-8 ---
$ws.Range("C2").Value = "NumberFormat: Date"
$ws.Range("C3").Value = "NumberFormat: Date"
$ws.Range("C4").Value = "Diebold 100, Aslam ASINH"
$ws.Range("C7").Value = "-"
$ws.Range("C5").Value = "Auto:market days will be counted based on data for each year, except for end year data using marketDaysYearEnd;Manual: all years using manualMarketDays"
$ws.Range("C6").Value = "Fill if using Manual"
$ws.Range("C8").Value = "Fill the probable market days for end year data / current year data"
$ws.Range("C9").Value = "Auto:decided by system using AIC;<integer>:fill if you know what the lag_order you want to use, example: 4"
$ws.Range("C10").Value = "Auto:default 10;<integer>:fill if you know what the forecast_horizon you want to use, example: 10"
$ws.Range("C11").Value = "Auto:default 200;<integer>:fill if you know what the rolling_window you want to use, example: 200"

# --- B values for the new rows (all reuse the existing "Auto" string) ---
$ws.Range("B9").Value = "Auto"
$ws.Range("B10").Value = "Auto"
$ws.Range("B11").Value = "Auto"

# Match style used for A2..A8 labels (style index 2) and B4/B5 dropdown values (style index 5)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9:A11").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B9:B11").PasteSpecial(-4122) | Out-Null

# Header cell C1 should match the header style used by A1/B1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update existing dateTo value (B3): 44196 -> 44195 (2020-12-31 -> 2020-12-30) ---
$ws.Range("B3").Value2 = 44195

# --- Selection moves to B3 ---
$ws.Range("B3").Select() | Out-Null
